# Regenerate s_vals data to filter save games.
# Updates columns B-E (and derived sum column G = B+C+D+E) for rows 2-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    3  = @(1.505614041169197, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    4  = @(0.1554434735375247, 0.05231270169004087, 16.98373111632243, 0.4998867070740569, 17.69137399862405)
    5  = @(0.7287194209349384, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 3.594575437922795)
    6  = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    7  = @(1.505614041169197, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    8  = @(0.7287194209349384, 1.65323645889881,  0.1529057820181812, 0.4998867070740569, 3.034748368925986)
    9  = @(3.182878228561681, 1.65323645889881,  3.082599426703578,  0.4998867070740569, 8.418600821238126)
    10 = @(0.1554434735375247, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 4.075514443323626)
    11 = @(3.182878228561681, 1.65323645889881,  0.7127328510149897, 6.48142807727062,   12.0302756157461)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G
}
